$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.654.95'
$ws.Range("E2").Value = '  +4.17%  '
$ws.Range("D3").Value = '1.797.38'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.36'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5298'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3766'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07533'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.121'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.45%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.64%  '
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.197'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.95%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.450'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.92%  '
$ws.Range("D16").Value = '1.790.95'
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.45'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.22%  '
$ws.Range("E18").Value = '  +0.67%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06444'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9998'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  +2.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.927'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").Value = '28.655.47'
$ws.Range("E23").Value = '  +4.02%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.096'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.99'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.56'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.405'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.65%  '
$ws.Range("D29").Value = '1.997.91'
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.55'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.45%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.131'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1020'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.731'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.88%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.664'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2298'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06577'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02327'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.805'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.077'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.53%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6322'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.46%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.205'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +5.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.396'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.57'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5925'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.665'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.06%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '125.78'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.981'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.00%  '
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06935'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.81%  '
